# Revert "Revert "Merge branch 'develop' ... into develop""
# Adds a missing "day 1" on-call tally (column C) for several teachers
# whose row was missing that entry; the dependent SUM/formula cells
# recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Month1")

$rows = @(9, 11, 12, 14, 28, 33, 47, 48, 49, 52, 56, 65)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 1
}

$wb.Application.CalculateFull()
